$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-149 down to 92-150
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with data (copy of old row 91's pattern,
# but with updated Fecha/Volumen/Precio values as per the new data point)
$ws.Cells.Item(91, 1).Value = 10
$ws.Cells.Item(91, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(91, 3).Value = "La Araucanía"
$ws.Cells.Item(91, 4).Value = 44719
$ws.Cells.Item(91, 4).NumberFormat = $ws.Cells.Item(92, 4).NumberFormat
$ws.Cells.Item(91, 5).Value = 9
$ws.Cells.Item(91, 6).Value = 100112012
$ws.Cells.Item(91, 7).Value = "Espinaca"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 30
$ws.Cells.Item(91, 11).Value = 12000
$ws.Cells.Item(91, 12).Value = 12000
$ws.Cells.Item(91, 13).Value = 12000
$ws.Cells.Item(91, 14).Value = "$/docena de atados"
$ws.Cells.Item(91, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(91, 16).Value = 4000
$ws.Cells.Item(91, 17).Value = 3
$ws.Cells.Item(91, 18).Value = "Hortaliza"
